$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "admin" block (rows 10-16): F changes from "U" to "VV" and a new
# rev/link-type-dkit / tabel-udah-jadi / fungsi-selesai tracking column
# gets filled in on columns G/H (and I on row 11). ---

# Row 10 - GET /admin
$ws.Range("F10").Value = "VV"
$ws.Range("G10").Value = "rev"
$ws.Range("H10").Value = "link type dkit"

# Row 11 - GET /admin/course
$ws.Range("F11").Value = "VV"
$ws.Range("G11").Value = "rev "
$ws.Range("H11").Value = "tabel udah jadi"
$ws.Range("I11").Value = "fungsi selesai"

# Row 12 - GET /admin/course/add
$ws.Range("F12").Value = "VV"
$ws.Range("G12").Value = "rev"
$ws.Range("H12").Value = "fungsi selesai"

# Row 13 - POST /admin/course/add
$ws.Range("F13").Value = "VV"
$ws.Range("G13").Value = "VV"
$ws.Range("H13").Value = "fungsi selesai"

# Row 14 - GET /admin/course/:courseId/edit
$ws.Range("F14").Value = "VV"
$ws.Range("G14").Value = "rev"
$ws.Range("H14").Value = "fungsi selesai"

# Row 15 - POST /admin/course/:courseId/edit
$ws.Range("F15").Value = "VV"
$ws.Range("G15").Value = "VV"
$ws.Range("H15").Value = "fungsi selesai"

# Row 16 - GET /admin/course/:courseId/delete
$ws.Range("F16").Value = "VV"
$ws.Range("G16").Value = "VV"
$ws.Range("H16").Value = "fungsi selesai"

# --- new trailing note row ---
$ws.Range("E29").Value = "SEARCHBAR"
$ws.Range("H29").Value = "BLOM DIO"

# --- column widths for the newly-used H/I columns ---
$ws.Range("H1").ColumnWidth = 14.67
$ws.Range("I1").ColumnWidth = 13.67

# --- view state: scroll so row 7 is visible near the top, and leave the
# active selection on H25 ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H25").Select()
